# Auto-generated edit script: updates static price/profit snapshot values
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) to match
# the refreshed market-data values captured by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H33").Value = 330.30768
$ws.Range("I33").Value = 210.2
$ws.Range("J33").Value = 730.6667
$ws.Range("K33").Value = 210.2
$ws.Range("L33").Value = 730.6667
$ws.Range("M33").Value = 18.80000000000001
$ws.Range("N33").Value = -1188.6667
$ws.Range("H40").Value = 4467.5
$ws.Range("J40").Value = 4999.8
$ws.Range("L40").Value = 4999.8
$ws.Range("N40").Value = -5349.8
$ws.Range("H92").Value = 435.85715
$ws.Range("I92").Value = 438.5909
$ws.Range("J92").Value = 425.83334
$ws.Range("K92").Value = 438.5909
$ws.Range("L92").Value = 425.83334
$ws.Range("M92").Value = 809.4091000000001
$ws.Range("N92").Value = -2921.83334
$ws.Range("H121").Value = 4999.5
$ws.Range("J121").Value = 4999.5
$ws.Range("L121").Value = 14998.5
$ws.Range("N121").Value = -18492.5
$ws.Range("H131").Value = 4229
$ws.Range("I131").Value = 2108.4285
$ws.Range("J131").Value = 7197.8
$ws.Range("K131").Value = 6325.2855
$ws.Range("L131").Value = 21593.4
$ws.Range("M131").Value = -1285.2855
$ws.Range("N131").Value = -31673.4
$ws.Range("H132").Value = 3946.4
$ws.Range("I132").Value = 4057.4119
$ws.Range("J132").Value = 3317.3333
$ws.Range("K132").Value = 12172.2357
$ws.Range("L132").Value = 9951.999899999999
$ws.Range("M132").Value = -9642.235700000001
$ws.Range("N132").Value = -15011.9999
$ws.Range("H137").Value = 497970.25
$ws.Range("I137").Value = 939490.2
$ws.Range("J137").Value = 12298.3
$ws.Range("K137").Value = 2818470.6
$ws.Range("L137").Value = 36894.89999999999
$ws.Range("M137").Value = -2815920.6
$ws.Range("N137").Value = -41994.89999999999
$ws.Range("H138").Value = 6222.9434
$ws.Range("I138").Value = 1816.6666
$ws.Range("J138").Value = 6785.447
$ws.Range("K138").Value = 5449.9998
$ws.Range("L138").Value = 20356.341
$ws.Range("M138").Value = -309.9997999999996
$ws.Range("N138").Value = -30636.341
$ws.Range("H141").Value = 5639.864
$ws.Range("I141").Value = 5130.25
$ws.Range("K141").Value = 15390.75
$ws.Range("M141").Value = -10210.75

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 7839.4287
$ws.Range("I2").Value = 12939
$ws.Range("J2").Value = 3203.4546
$ws.Range("K2").Value = 12939
$ws.Range("L2").Value = 3203.4546
$ws.Range("M2").Value = -12826
$ws.Range("N2").Value = -3429.4546
$ws.Range("H32").Value = 2162.365
$ws.Range("I32").Value = 1970.5667
$ws.Range("J32").Value = 5998.3335
$ws.Range("K32").Value = 1970.5667
$ws.Range("L32").Value = 5998.3335
$ws.Range("M32").Value = -1683.5667
$ws.Range("N32").Value = -6572.3335
$ws.Range("H61").Value = 6567.5835
$ws.Range("I61").Value = 7462.4443
$ws.Range("J61").Value = 3883
$ws.Range("K61").Value = 7462.4443
$ws.Range("L61").Value = 3883
$ws.Range("M61").Value = -7250.4443
$ws.Range("N61").Value = -4307
$ws.Range("H116").Value = 7839.4287
$ws.Range("I116").Value = 12939
$ws.Range("J116").Value = 3203.4546
$ws.Range("K116").Value = 12939
$ws.Range("L116").Value = 3203.4546
$ws.Range("M116").Value = -10645
$ws.Range("N116").Value = -7791.4546
$ws.Range("H136").Value = 6567.5835
$ws.Range("I136").Value = 7462.4443
$ws.Range("J136").Value = 3883
$ws.Range("K136").Value = 22387.3329
$ws.Range("L136").Value = 11649
$ws.Range("M136").Value = -19837.3329
$ws.Range("N136").Value = -16749

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 7839.4287
$ws.Range("I3").Value = 12939
$ws.Range("J3").Value = 3203.4546
$ws.Range("K3").Value = 12939
$ws.Range("L3").Value = 3203.4546
$ws.Range("M3").Value = -12825
$ws.Range("N3").Value = -3431.4546
$ws.Range("H107").Value = 2179.423
$ws.Range("I107").Value = 1948.4762
$ws.Range("J107").Value = 3149.4
$ws.Range("K107").Value = 1948.4762
$ws.Range("L107").Value = 3149.4
$ws.Range("M107").Value = -28.47620000000006
$ws.Range("N107").Value = -6989.4

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2244.5103
$ws.Range("I31").Value = 1662.6111
$ws.Range("J31").Value = 2582.3872
$ws.Range("K31").Value = 1662.6111
$ws.Range("L31").Value = 2582.3872
$ws.Range("M31").Value = -1367.6111
$ws.Range("N31").Value = -3172.3872
$ws.Range("H34").Value = 2244.5103
$ws.Range("I34").Value = 1662.6111
$ws.Range("J34").Value = 2582.3872
$ws.Range("K34").Value = 1662.6111
$ws.Range("L34").Value = 2582.3872
$ws.Range("M34").Value = -1460.6111
$ws.Range("N34").Value = -2986.3872
$ws.Range("H58").Value = 5177.595
$ws.Range("I58").Value = 5059.4243
$ws.Range("K58").Value = 5059.4243
$ws.Range("M58").Value = -4856.4243
$ws.Range("H134").Value = 2983333.5
$ws.Range("I134").Value = 3480139
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 10440417
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -10437882
$ws.Range("N134").Value = -12570
$ws.Range("H136").Value = 5177.595
$ws.Range("I136").Value = 5059.4243
$ws.Range("K136").Value = 15178.2729
$ws.Range("M136").Value = -12628.2729

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H64").Value = 2348.1667
$ws.Range("I64").Value = 1863
$ws.Range("J64").Value = 2833.3333
$ws.Range("K64").Value = 5589
$ws.Range("L64").Value = 8499.999899999999
$ws.Range("M64").Value = -5319
$ws.Range("N64").Value = -9039.999899999999
$ws.Range("H67").Value = 2348.1667
$ws.Range("I67").Value = 1863
$ws.Range("J67").Value = 2833.3333
$ws.Range("K67").Value = 5589
$ws.Range("L67").Value = 8499.999899999999
$ws.Range("M67").Value = -4653
$ws.Range("N67").Value = -10371.9999
$ws.Range("H68").Value = 15158838
$ws.Range("I68").Value = 2874.0833
$ws.Range("J68").Value = 23819390
$ws.Range("K68").Value = 8622.249899999999
$ws.Range("L68").Value = 71458170
$ws.Range("M68").Value = -7811.249899999999
$ws.Range("N68").Value = -71459792
$ws.Range("H69").Value = 4600
$ws.Range("J69").Value = 4600
$ws.Range("L69").Value = 13800
$ws.Range("N69").Value = -15422
$ws.Range("H71").Value = 15158838
$ws.Range("I71").Value = 2874.0833
$ws.Range("J71").Value = 23819390
$ws.Range("K71").Value = 25866.7497
$ws.Range("L71").Value = 214374510
$ws.Range("M71").Value = -21810.7497
$ws.Range("N71").Value = -214382622
$ws.Range("H72").Value = 4600
$ws.Range("J72").Value = 4600
$ws.Range("L72").Value = 41400
$ws.Range("N72").Value = -49512
$ws.Range("H76").Value = 8966.799999999999
$ws.Range("I76").Value = 8966.799999999999
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 26900.4
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -26517.4
$ws.Range("N76").Value = $null
$ws.Range("H79").Value = 8966.799999999999
$ws.Range("I79").Value = 8966.799999999999
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 26900.4
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -25574.4
$ws.Range("N79").Value = $null
$ws.Range("H80").Value = 110188.664
$ws.Range("I80").Value = 2815.5
$ws.Range("J80").Value = 181770.78
$ws.Range("K80").Value = 8446.5
$ws.Range("L80").Value = 545312.34
$ws.Range("M80").Value = -7510.5
$ws.Range("N80").Value = -547184.34
$ws.Range("H83").Value = 110188.664
$ws.Range("I83").Value = 2815.5
$ws.Range("J83").Value = 181770.78
$ws.Range("K83").Value = 25339.5
$ws.Range("L83").Value = 1635937.02
$ws.Range("M83").Value = -20659.5
$ws.Range("N83").Value = -1645297.02
$ws.Range("H131").Value = 21279998
$ws.Range("I131").Value = 76929496
$ws.Range("J131").Value = 2249.0881
$ws.Range("K131").Value = 230788488
$ws.Range("L131").Value = 6747.2643
$ws.Range("M131").Value = -230783448
$ws.Range("N131").Value = -16827.2643

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 4648.3516
$ws.Range("I70").Value = 4334.9287
$ws.Range("J70").Value = 5623.4443
$ws.Range("K70").Value = 4334.9287
$ws.Range("L70").Value = 5623.4443
$ws.Range("M70").Value = -4064.9287
$ws.Range("N70").Value = -6163.4443
$ws.Range("H73").Value = 4648.3516
$ws.Range("I73").Value = 4334.9287
$ws.Range("J73").Value = 5623.4443
$ws.Range("K73").Value = 4334.9287
$ws.Range("L73").Value = 5623.4443
$ws.Range("M73").Value = -3398.9287
$ws.Range("N73").Value = -7495.4443
$ws.Range("H102").Value = 33970.75
$ws.Range("I102").Value = 30664.9
$ws.Range("J102").Value = 50500
$ws.Range("K102").Value = 30664.9
$ws.Range("L102").Value = 50500
$ws.Range("M102").Value = -29042.9
$ws.Range("N102").Value = -53744
$ws.Range("H122").Value = 16926.2
$ws.Range("I122").Value = 21966.143
$ws.Range("J122").Value = 5166.3335
$ws.Range("K122").Value = 65898.429
$ws.Range("L122").Value = 15499.0005
$ws.Range("M122").Value = -63448.429
$ws.Range("N122").Value = -20399.0005

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 2946.5
$ws.Range("I46").Value = 1187.3334
$ws.Range("J46").Value = 4164.385
$ws.Range("K46").Value = 1187.3334
$ws.Range("L46").Value = 4164.385
$ws.Range("M46").Value = -999.3334
$ws.Range("N46").Value = -4540.385
$ws.Range("H68").Value = 6328.357
$ws.Range("I68").Value = 9501
$ws.Range("J68").Value = 5799.5835
$ws.Range("K68").Value = 9501
$ws.Range("L68").Value = 5799.5835
$ws.Range("M68").Value = -8752
$ws.Range("N68").Value = -7297.5835
$ws.Range("H71").Value = 6328.357
$ws.Range("I71").Value = 9501
$ws.Range("J71").Value = 5799.5835
$ws.Range("K71").Value = 47505
$ws.Range("L71").Value = 28997.9175
$ws.Range("M71").Value = -43761
$ws.Range("N71").Value = -36485.9175
$ws.Range("H136").Value = 5924.4165
$ws.Range("I136").Value = 2824.125
$ws.Range("J136").Value = 12125
$ws.Range("K136").Value = 8472.375
$ws.Range("L136").Value = 36375
$ws.Range("M136").Value = -5922.375
$ws.Range("N136").Value = -41475
